$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    if ($text -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

Set-TextValue 'D2' '30.312.77'
Set-TextValue 'E2' '  +0.03%  '
Set-TextValue 'D3' '1.932.17'
Set-TextValue 'E3' '  +0.10%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'D5' '0.7488'
Set-TextValue 'E5' '  +4.52%  '
Set-TextValue 'D6' '243.24'
Set-TextValue 'E6' '  -2.30%  '
Set-TextValue 'E7' '  +0.08%  '
Set-TextValue 'D8' '27.61'
Set-TextValue 'E8' '  -0.27%  '
Set-TextValue 'D9' '0.3170'
Set-TextValue 'E9' '  -1.24%  '
Set-TextValue 'D10' '0.07112'
Set-TextValue 'E10' '  +0.10%  '
Set-TextValue 'D11' '0.08058'
Set-TextValue 'E11' '  +0.73%  '
Set-TextValue 'D12' '0.7783'
Set-TextValue 'E12' '  -1.67%  '
Set-TextValue 'D13' '1.890.72'
Set-TextValue 'E13' '  -2.08%  '
Set-TextValue 'D14' '5.395'
Set-TextValue 'D15' '93.07'
Set-TextValue 'E15' '  -1.93%  '
Set-TextValue 'D16' '14.55'
Set-TextValue 'E16' '  -0.78%  '
Set-TextValue 'D17' '30.311.91'
Set-TextValue 'E17' '  -0.01%  '
Set-TextValue 'D18' '6.019'
Set-TextValue 'E18' '  +4.25%  '
Set-TextValue 'D19' '251.86'
Set-TextValue 'E19' '  -2.18%  '
Set-TextValue 'D20' '0.000007939'
Set-TextValue 'E20' '  -1.57%  '
Set-TextValue 'B21' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C21' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D21' '2.167.66'
Set-TextValue 'E21' '  -0.55%  '
Set-TextValue 'B22' 'Dai'
Set-TextValue 'C22' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D22' '1.001'
Set-TextValue 'E22' '  +0.04%  '
Set-TextValue 'E23' '  +0.11%  '
Set-TextValue 'D24' '6.679'
Set-TextValue 'E24' '  -2.32%  '
Set-TextValue 'D25' '9.550'
Set-TextValue 'E25' '  +0.00%  '
Set-TextValue 'D26' '165.18'
Set-TextValue 'E26' '  +0.21%  '
Set-TextValue 'D27' '19.07'
Set-TextValue 'E27' '  -0.26%  '
Set-TextValue 'D28' '0.1296'
Set-TextValue 'E28' '  +1.91%  '
Set-TextValue 'D29' '2.184'
Set-TextValue 'E29' '  -3.98%  '
Set-TextValue 'D30' '1.369'
Set-TextValue 'E30' '  +1.10%  '
Set-TextValue 'D31' '1.560'
Set-TextValue 'E31' '  +2.22%  '
Set-TextValue 'D32' '4.410'
Set-TextValue 'E32' '  +0.22%  '
Set-TextValue 'D33' '4.141'
Set-TextValue 'E33' '  -0.07%  '
Set-TextValue 'D34' '0.05228'
Set-TextValue 'E34' '  +1.60%  '
Set-TextValue 'D35' '1.320'
Set-TextValue 'E35' '  +4.10%  '
Set-TextValue 'D36' '0.7570'
Set-TextValue 'E36' '  +1.66%  '
Set-TextValue 'D37' '2.782'
Set-TextValue 'E37' '  +0.55%  '
Set-TextValue 'D38' '0.01947'
Set-TextValue 'E38' '  -0.92%  '
Set-TextValue 'D39' '2.798'
Set-TextValue 'E39' '  +0.02%  '
Set-TextValue 'D40' '6.513'
Set-TextValue 'E40' '  +2.37%  '
Set-TextValue 'D41' '78.16'
Set-TextValue 'E41' '  -0.71%  '
Set-TextValue 'D42' '0.4522'
Set-TextValue 'E42' '  +0.13%  '
Set-TextValue 'D43' '1.977'
Set-TextValue 'E43' '  -1.07%  '
Set-TextValue 'D44' '0.8421'
Set-TextValue 'E44' '  -0.72%  '
Set-TextValue 'D46' '7.690'
Set-TextValue 'E46' '  +3.33%  '
Set-TextValue 'D47' '9.958'
Set-TextValue 'E47' '  +1.75%  '
Set-TextValue 'D48' '101.58'
Set-TextValue 'E48' '  +1.02%  '
Set-TextValue 'D49' '37.97'
Set-TextValue 'E49' '  +3.35%  '
Set-TextValue 'D50' '0.1230'
Set-TextValue 'E50' '  +7.07%  '
Set-TextValue 'D51' '960.12'
Set-TextValue 'E51' '  +1.86%  '
